$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$ws.Range("D5").Value = "dsfjds"
$ws.Range("D5").Select()
